$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2023-10-08 10:42:36", "hatespeech", "setting1", 10, 0.6832539858573135)
    ,@("2023-10-08 10:42:36", "hatespeech", "def", 10, 0.6824274371121277)
    ,@("2023-10-08 10:42:36", "hatespeech", "setting5", 20, 0.5430190088399398)
    ,@("2023-10-08 10:42:36", "hatespeech", "setting2", 10, 0.6888464370033944)
    ,@("2023-10-08 10:42:36", "hatespeech", "setting4", 10, 0.685377026532934)
    ,@("2023-10-08 10:42:37", "hatespeech", "setting5", 10, 0.6866927584874289)
    ,@("2023-10-08 10:42:37", "hatespeech", "setting6", 20, 0.5372232041800609)
    ,@("2023-10-08 10:42:37", "hatespeech", "def", 20, 0.5370654831575581)
    ,@("2023-10-08 10:42:37", "hatespeech", "setting1", 20, 0.529906608328365)
    ,@("2023-10-08 10:42:37", "hatespeech", "setting2", 20, 0.5463945139545814)
    ,@("2023-10-08 10:42:37", "hatespeech", "setting3", 10, 0.6777993948529082)
    ,@("2023-10-08 10:42:37", "hatespeech", "setting6", 10, 0.6832857918439053)
    ,@("2023-10-08 10:42:37", "hatespeech", "setting3", 20, 0.5304085584820892)
    ,@("2023-10-08 10:42:38", "hatespeech", "setting1", 30, 0.4372067598742038)
    ,@("2023-10-08 10:42:38", "hatespeech", "setting4", 20, 0.5437831482888089)
    ,@("2023-10-08 10:42:38", "hatespeech", "setting2", 30, 0.4512105366595455)
    ,@("2023-10-08 10:42:38", "hatespeech", "setting5", 30, 0.4561276175217827)
    ,@("2023-10-08 10:42:38", "hatespeech", "setting6", 30, 0.448608784273391)
    ,@("2023-10-08 10:42:38", "hatespeech", "def", 30, 0.435017901375769)
    ,@("2023-10-08 10:42:39", "hatespeech", "setting1", 40, 0.3608176208054734)
    ,@("2023-10-08 10:42:39", "hatespeech", "setting3", 30, 0.4429020612955527)
    ,@("2023-10-08 10:42:39", "hatespeech", "setting2", 40, 0.3748490918771323)
    ,@("2023-10-08 10:42:39", "hatespeech", "setting4", 30, 0.448449170617598)
    ,@("2023-10-08 10:42:39", "hatespeech", "setting3", 40, 0.3726524685589809)
    ,@("2023-10-08 10:42:39", "hatespeech", "def", 40, 0.3625392394626541)
    ,@("2023-10-08 10:42:39", "hatespeech", "setting6", 40, 0.3734577247715302)
    ,@("2023-10-08 10:42:40", "hatespeech", "setting4", 40, 0.3837302197909441)
    ,@("2023-10-08 10:42:40", "hatespeech", "setting5", 40, 0.3683956783844847)
)

$startRow = 149
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}